$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 316, shifting existing rows 316-414 down to 317-415.
$ws.Rows.Item(316).EntireRow.Insert()

# Populate the new row 316 with the new data record.
$ws.Range("A316").Value = 6
$ws.Range("B316").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C316").Value = "Metropolitana"
$ws.Range("D316").Value = 44559
$ws.Range("D316").NumberFormat = $ws.Range("D317").NumberFormat
$ws.Range("E316").Value = 13
$ws.Range("F316").Value = 100112044
$ws.Range("G316").Value = "Perejil"
$ws.Range("H316").Value = "Sin especificar"
$ws.Range("I316").Value = "Primera"
$ws.Range("J316").Value = 230
$ws.Range("K316").Value = 9000
$ws.Range("L316").Value = 10000
$ws.Range("M316").Value = 9391
$ws.Range("N316").Value = "$/docena de atados"
$ws.Range("O316").Value = "Región Metropolitana"
$ws.Range("P316").Value = 3130
$ws.Range("Q316").Value = 3
$ws.Range("R316").Value = "Hortaliza"
